# Actualización desde MV -datos-
# Adds the newly reported daily auction rows (10-08-2021 .. 02-09-2021)
# to the bottom of the existing table, following the same layout as the
# pre-existing rows (A: date label, B-F: amounts, G: rate).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newRows = @(
    @("10-08-2021", 10000, 27000, 10000, 10000, 0,    1.5),
    @("11-08-2021", 10000, $null, 0,     $null, $null, $null),
    @("12-08-2021", 10000, 30000, 10000, 10000, 0,    1.55),
    @("17-08-2021", 10000, 35000, 10000, 8000,  2000, 1.6),
    @("18-08-2021", 10000, 31000, 10000, 10000, 0,    1.6),
    @("19-08-2021", 10000, 22000, 10000, 8000,  2000, 1.65),
    @("24-08-2021", 10000, 37000, 10000, 10000, 0,    1.65),
    @("25-08-2021", 10000, 27000, 10000, 10000, 0,    1.65),
    @("26-08-2021", 10000, 21000, 10000, 10000, 0,    1.64),
    @("01-09-2021", 10000, 19000, 10000, 10000, 0,    2.25),
    @("02-09-2021", 10000, $null, 0,     $null, $null, $null)
)

$startRow = 8
$endRow = $startRow + $newRows.Count - 1

# Force column A on these rows to be treated as plain text so labels like
# "10-08-2021" / "01-09-2021" aren't silently reinterpreted as date serials
# (day <= 12, so they're ambiguous/date-like strings).
$ws.Range("A$startRow`:A$endRow").NumberFormat = "@"

for ($i = 0; $i -lt $newRows.Count; $i++) {
    $r = $startRow + $i
    $data = $newRows[$i]

    $ws.Cells.Item($r, 1).Value = $data[0]

    if ($null -ne $data[1]) { $ws.Cells.Item($r, 2).Value = $data[1] }
    if ($null -ne $data[2]) { $ws.Cells.Item($r, 3).Value = $data[2] }
    if ($null -ne $data[3]) { $ws.Cells.Item($r, 4).Value = $data[3] }
    if ($null -ne $data[4]) { $ws.Cells.Item($r, 5).Value = $data[4] }
    if ($null -ne $data[5]) { $ws.Cells.Item($r, 6).Value = $data[5] }
    if ($null -ne $data[6]) { $ws.Cells.Item($r, 7).Value = $data[6] }
}

# Restore the default ("Normal") style on the new A-column cells now that
# the values are locked in as text, so no residual number-format styling
# is left on the cells themselves.
$ws.Range("A$startRow`:A$endRow").Style = "Normal"
